# Insert a new weekly price record as row 190 on the "Ajo" (Vega Monumental
# Concepción) sheet, pushing the existing rows 190-294 down to 191-295.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 190:294 down by inserting a new blank row at 190.
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A190").Value = 11
$ws.Range("B190").Value = "Vega Monumental Concepción"
$ws.Range("C190").Value = "Bíobío"
$ws.Range("D190").Value = "2023-06-20"
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = 100112003
$ws.Range("G190").Value = "Ajo"
$ws.Range("H190").Value = "Chino"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 220
$ws.Range("K190").Value = 15000
$ws.Range("L190").Value = 16000
$ws.Range("M190").Value = 15455
$ws.Range("N190").Value = "$/caja 10 kilos"
$ws.Range("O190").Value = "China"
$ws.Range("P190").Value = 1546
$ws.Range("Q190").Value = 10
$ws.Range("R190").Value = "Hortaliza"
